# lines_states.xlsx / Sheet1 edit
#
# Inserts two new data rows ("line7" and "line8") right after the existing
# "line6" row (i.e. before the old row 8), pushing the "extr1".."extr8"
# rows down by two. The "extr*" rows keep exactly the same from_bus /
# to_bus / in_service values they had before - only their row position and
# their auto-numbered index in column A change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the two new rows, shifting old rows 8:15 ("extr1".."extr8")
# down to 10:17.
$ws.Rows("8:9").Insert() | Out-Null

# The insert leaves the new rows 8:9 without the bordered/bold/centered
# style used by every other data row in column A (style index 1 in
# styles.xml). Copy formatting from an existing data row (row 2) onto the
# new rows so they match the rest of the table.
$ws.Range("A2:E2").Copy() | Out-Null
$ws.Range("A8:E9").PasteSpecial(-4122) | Out-Null

# New row for "line7"
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

# New row for "line8"
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $false

# The "extr1".."extr8" rows (now at 10:17) keep their original data, but
# the running index in column A needs to be renumbered to account for the
# two rows inserted above them.
for ($i = 10; $i -le 17; $i++) {
  $ws.Range("A$i").Value = $i - 2
}
